$d = $word.ActiveDocument

# Locate the title run "FETAL DIAGNOSIS ENHANCEMENT TOOL" (search only, no
# replace yet) so we can compute the exact character offset that sits
# between "ENHANCEMENT " and "TOOL".
$titleRange = $d.Content
$titleRange.Find.Execute("FETAL DIAGNOSIS ENHANCEMENT TOOL")

$prefix = "FETAL DIAGNOSIS ENHANCEMENT "
$bookmarkPos = $titleRange.Start + $prefix.Length

# Re-seat the (now orphaned) "_GoBack" bookmark right after "ENHANCEMENT "
# and before "TOOL" -- adding a bookmark with an existing name moves it,
# which also removes it from its old location at the end of the document.
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Now trim the trailing " TOOL" word, operating only on the text after the
# bookmark so the edit doesn't span (and delete) the bookmark we just set.
$tailRange = $d.Range($bookmarkPos, $d.Content.End)
$tailRange.Find.Execute("TOOL", $true, $false, $false, $false, $false, $true,
                         1, $false, "", 2)
